$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix two datetime values that were off by 1 hour (A6, A24) ---
$ws.Range("A6").Value = 38443.45833333334
$ws.Range("A24").Value = 38991.45833333334

# --- Revise the M2 figures for rows 218-220 (open/high/low/close) ---
$ws.Range("C218").Value = 603042000000
$ws.Range("D218").Value = 603042000000
$ws.Range("E218").Value = 603042000000
$ws.Range("F218").Value = 603042000000

$ws.Range("C219").Value = 603647900000
$ws.Range("D219").Value = 603647900000
$ws.Range("E219").Value = 603647900000
$ws.Range("F219").Value = 603647900000

$ws.Range("C220").Value = 612183900000
$ws.Range("D220").Value = 612183900000
$ws.Range("E220").Value = 612183900000
$ws.Range("F220").Value = 612183900000

# --- Append a new row 224 with the latest data point ---
$ws.Range("A224").Value = 45078.41666666666
$ws.Range("B224").Value = "ECONOMICS:ROM2"
$ws.Range("C224").Value = 624519300000
$ws.Range("D224").Value = 624519300000
$ws.Range("E224").Value = 624519300000
$ws.Range("F224").Value = 624519300000
$ws.Range("G224").Value = 0

# Match the formatting used by the rest of the table (row 223) so the new
# row's cell styles - in particular the datetime number format in column A -
# line up with the existing data rows.
$ws.Range("A223:G223").Copy()
$ws.Range("A224").PasteSpecial(-4122)

# Re-apply the values, since PasteSpecial(xlPasteFormats) only copies
# formatting, but do it defensively in case paste also touched values.
$ws.Range("A224").Value = 45078.41666666666
$ws.Range("B224").Value = "ECONOMICS:ROM2"
$ws.Range("C224").Value = 624519300000
$ws.Range("D224").Value = 624519300000
$ws.Range("E224").Value = 624519300000
$ws.Range("F224").Value = 624519300000
$ws.Range("G224").Value = 0
